$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.015.11"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.897.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5019"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3919"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09338"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.130"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.99"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.358"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.74"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.897.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.296"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.44"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06575"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.211"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.066.91"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.320"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.627"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.119.81"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.86"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.50"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.51"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.079"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1065"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.605"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.623"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.606"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06612"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02424"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.227"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2173"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.279"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.43%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.993"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6343"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.40"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5964"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.713"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.273"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.034"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.28"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.176"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.46%  "
